$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rodada 1")

# Mandante_Pontos (column D) and Visitante_Pontos (column F) for rows 2-17
$values = @{
    2  = @{ D = 38.5;              F = 39.66 }
    3  = @{ D = 41.6;              F = 64.96 }
    4  = @{ D = 59.76;             F = 45.86 }
    5  = @{ D = 39.66;             F = 58.4 }
    6  = @{ D = 34.36;             F = 71.36 }
    7  = @{ D = 69.56;             F = 53.06 }
    8  = @{ D = 69.76;             F = 45.46 }
    9  = @{ D = 30.6;              F = 59.56 }
    10 = @{ D = 50.6;              F = 73.66 }
    11 = @{ D = 49.16;             F = 57.25 }
    12 = @{ D = 30.06;             F = 54.36 }
    13 = @{ D = 47.2;              F = 51.05 }
    14 = @{ D = 45.3;              F = 57.76 }
    15 = @{ D = 69.26;             F = 43.1 }
    16 = @{ D = 34.76;             F = 56.65 }
    17 = @{ D = 59.86;             F = 50.76 }
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row].D
    $ws.Range("F$row").Value = $values[$row].F
}
